# beril_model.xlsx update — "collection examples inc DH"
#
# - Adds a NamedThingCollection sheet (material_entities, processes) in front,
#   reusing the old "NamedThing" sheet's slot and recreating NamedThing itself
#   (id, name, description, type_value) right after it.
# - Every concrete entity sheet (NamedThing, MaterialEntity, InformationArtifact,
#   Process) gains a trailing `type_value` column.
# - MaterialEntity drops its `observations` column.
# - Two new physical-sample sheets are added after MaterialEntity: SoilSample and
#   DnaExtract, both (mass_g, color, id, name, description, type_value) with a
#   RED/GREEN/BLUE dropdown validation on `color`.
# - Person and PersonCollection (and the vital_status validation that lived on
#   Person) are removed entirely.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. NamedThing -> becomes the new NamedThingCollection sheet; a fresh
#    NamedThing sheet is inserted right after it with the old column set plus
#    type_value.
# ---------------------------------------------------------------------------
$namedThingCollection = $wb.Worksheets.Item("NamedThing")
$namedThingCollection.Name = "NamedThingCollection"
$namedThingCollection.Range("A1").Value = "material_entities"
$namedThingCollection.Range("B1").Value = "processes"
$namedThingCollection.Range("C1").ClearContents()

$namedThing = $wb.Worksheets.Add($null, $namedThingCollection)
$namedThing.Name = "NamedThing"
$namedThing.Range("A1").Value = "id"
$namedThing.Range("B1").Value = "name"
$namedThing.Range("C1").Value = "description"
$namedThing.Range("D1").Value = "type_value"

# ---------------------------------------------------------------------------
# 2. Observation / NonProcess are untouched.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3. MaterialEntity loses `observations`, gains `type_value`.
# ---------------------------------------------------------------------------
$materialEntity = $wb.Worksheets.Item("MaterialEntity")
$materialEntity.Range("A1").Value = "id"
$materialEntity.Range("B1").Value = "name"
$materialEntity.Range("C1").Value = "description"
$materialEntity.Range("D1").Value = "type_value"

# ---------------------------------------------------------------------------
# 4. Two brand-new sample sheets, inserted right after MaterialEntity.
# ---------------------------------------------------------------------------
$soilSample = $wb.Worksheets.Add($null, $materialEntity)
$soilSample.Name = "SoilSample"
$soilSample.Range("A1").Value = "mass_g"
$soilSample.Range("B1").Value = "color"
$soilSample.Range("C1").Value = "id"
$soilSample.Range("D1").Value = "name"
$soilSample.Range("E1").Value = "description"
$soilSample.Range("F1").Value = "type_value"
$soilSample.Range("B2:B1048576").Validation.Add(3, 1, 1, '"RED,GREEN,BLUE"')

$dnaExtract = $wb.Worksheets.Add($null, $soilSample)
$dnaExtract.Name = "DnaExtract"
$dnaExtract.Range("A1").Value = "mass_g"
$dnaExtract.Range("B1").Value = "color"
$dnaExtract.Range("C1").Value = "id"
$dnaExtract.Range("D1").Value = "name"
$dnaExtract.Range("E1").Value = "description"
$dnaExtract.Range("F1").Value = "type_value"
$dnaExtract.Range("B2:B1048576").Validation.Add(3, 1, 1, '"RED,GREEN,BLUE"')

# ---------------------------------------------------------------------------
# 5. InformationArtifact gains `type_value`.
# ---------------------------------------------------------------------------
$infoArtifact = $wb.Worksheets.Item("InformationArtifact")
$infoArtifact.Range("A1").Value = "size_in_bytes"
$infoArtifact.Range("B1").Value = "md5"
$infoArtifact.Range("C1").Value = "url"
$infoArtifact.Range("D1").Value = "id"
$infoArtifact.Range("E1").Value = "name"
$infoArtifact.Range("F1").Value = "description"
$infoArtifact.Range("G1").Value = "type_value"

# ---------------------------------------------------------------------------
# 6. Process gains `type_value`.
# ---------------------------------------------------------------------------
$process = $wb.Worksheets.Item("Process")
$process.Range("A1").Value = "inputs"
$process.Range("B1").Value = "outputs"
$process.Range("C1").Value = "id"
$process.Range("D1").Value = "name"
$process.Range("E1").Value = "description"
$process.Range("F1").Value = "type_value"

# ---------------------------------------------------------------------------
# 7. Person / PersonCollection are dropped (and with them the vital_status
#    validation that lived on Person).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Person").Delete() | Out-Null
$wb.Worksheets.Item("PersonCollection").Delete() | Out-Null

# Keep the first sheet active/selected, same as before the edit.
$namedThingCollection.Activate()
$namedThingCollection.Range("A1").Select() | Out-Null
